$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1662.909
$ws.Range("J112").Value = 2161.5
$ws.Range("L112").Value = 6484.5
$ws.Range("N112").Value = -8700.5

$ws.Range("H125").Value = 83342216
$ws.Range("J125").Value = 100010550
$ws.Range("L125").Value = 900094950
$ws.Range("N125").Value = -900099870

$ws.Range("H129").Value = 1065.8701
$ws.Range("J129").Value = 1065.8701
$ws.Range("L129").Value = 3197.6103
$ws.Range("N129").Value = -13197.6103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18163.25
$ws.Range("I32").Value = 19279.982
$ws.Range("K32").Value = 19279.982
$ws.Range("M32").Value = -18992.982

$ws.Range("H37").Value = 28971.428
$ws.Range("J37").Value = 28971.428
$ws.Range("L37").Value = 28971.428
$ws.Range("N37").Value = -29517.428

$ws.Range("H61").Value = 8425.093000000001
$ws.Range("I61").Value = 5263.2905
$ws.Range("J61").Value = 16593.084
$ws.Range("K61").Value = 5263.2905
$ws.Range("L61").Value = 16593.084
$ws.Range("M61").Value = -5051.2905
$ws.Range("N61").Value = -17017.084

$ws.Range("H74").Value = 208699.17
$ws.Range("I74").Value = 272239.38
$ws.Range("J74").Value = 2193.5
$ws.Range("K74").Value = 272239.38
$ws.Range("L74").Value = 2193.5
$ws.Range("M74").Value = -271365.38
$ws.Range("N74").Value = -3941.5

$ws.Range("H77").Value = 208699.17
$ws.Range("I77").Value = 272239.38
$ws.Range("J77").Value = 2193.5
$ws.Range("K77").Value = 1361196.9
$ws.Range("L77").Value = 10967.5
$ws.Range("M77").Value = -1356828.9
$ws.Range("N77").Value = -19703.5

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H132").Value = 3731.476
$ws.Range("I132").Value = 3279.0833
$ws.Range("K132").Value = 9837.249899999999
$ws.Range("M132").Value = -7307.249899999999

$ws.Range("H136").Value = 8425.093000000001
$ws.Range("I136").Value = 5263.2905
$ws.Range("J136").Value = 16593.084
$ws.Range("K136").Value = 15789.8715
$ws.Range("L136").Value = 49779.25199999999
$ws.Range("M136").Value = -13239.8715
$ws.Range("N136").Value = -54879.25199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11667083
$ws.Range("I6").Value = 11667083
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 11667083
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -11666970
$ws.Range("N6").ClearContents()

$ws.Range("H32").Value = 3802.2
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 3011
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 3011
$ws.Range("M32").Value = -3684
$ws.Range("N32").Value = -3643

$ws.Range("H35").Value = 841.6667
$ws.Range("I35").Value = 841.6667
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 841.6667
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -547.6667
$ws.Range("N35").ClearContents()

$ws.Range("H58").Value = 2843355
$ws.Range("I58").Value = 6494717
$ws.Range("J58").Value = 3406.5557
$ws.Range("K58").Value = 6494717
$ws.Range("L58").Value = 3406.5557
$ws.Range("M58").Value = -6494514
$ws.Range("N58").Value = -3812.5557

$ws.Range("H132").Value = 2692.9092
$ws.Range("I132").Value = 2492.8572
$ws.Range("K132").Value = 7478.571599999999
$ws.Range("M132").Value = -4948.571599999999

$ws.Range("H136").Value = 2843355
$ws.Range("I136").Value = 6494717
$ws.Range("J136").Value = 3406.5557
$ws.Range("K136").Value = 19484151
$ws.Range("L136").Value = 10219.6671
$ws.Range("M136").Value = -19481601
$ws.Range("N136").Value = -15319.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21869.762
$ws.Range("I131").Value = 1859.091
$ws.Range("J131").Value = 28158.828
$ws.Range("K131").Value = 5577.272999999999
$ws.Range("L131").Value = 84476.484
$ws.Range("M131").Value = -537.2729999999992
$ws.Range("N131").Value = -94556.484

$ws.Range("H132").Value = 2786.2856
$ws.Range("I132").Value = 5502
$ws.Range("K132").Value = 49518
$ws.Range("M132").Value = -46988

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 18000
$ws.Range("J69").Value = 18000
$ws.Range("L69").Value = 18000
$ws.Range("N69").Value = -19498

$ws.Range("H72").Value = 18000
$ws.Range("J72").Value = 18000
$ws.Range("L72").Value = 54000
$ws.Range("N72").Value = -61488

$ws.Range("H132").Value = 10820.143
$ws.Range("I132").Value = 3746.9
$ws.Range("J132").Value = 28503.25
$ws.Range("K132").Value = 11240.7
$ws.Range("L132").Value = 85509.75
$ws.Range("M132").Value = -8710.700000000001
$ws.Range("N132").Value = -90569.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5777.615
$ws.Range("I7").Value = 5838
$ws.Range("J7").Value = 5681
$ws.Range("K7").Value = 5838
$ws.Range("L7").Value = 5681
$ws.Range("M7").Value = -5726
$ws.Range("N7").Value = -5905

$ws.Range("H40").Value = 2768.5
$ws.Range("I40").Value = 2380.65
$ws.Range("J40").Value = 3544.2
$ws.Range("K40").Value = 2380.65
$ws.Range("L40").Value = 3544.2
$ws.Range("M40").Value = -2244.65
$ws.Range("N40").Value = -3816.2

$ws.Range("H126").Value = 5777.615
$ws.Range("I126").Value = 5838
$ws.Range("J126").Value = 5681
$ws.Range("K126").Value = 17514
$ws.Range("L126").Value = 17043
$ws.Range("M126").Value = -15044
$ws.Range("N126").Value = -21983

$ws.Range("H132").Value = 4290.4443
$ws.Range("I132").Value = 4309.75
$ws.Range("J132").Value = 4262.364
$ws.Range("K132").Value = 12929.25
$ws.Range("L132").Value = 12787.092
$ws.Range("M132").Value = -10399.25
$ws.Range("N132").Value = -17847.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3919.5454
$ws.Range("J62").Value = 3889.375
$ws.Range("L62").Value = 3889.375
$ws.Range("N62").Value = -5137.375

$ws.Range("H64").Value = 40114
$ws.Range("J64").Value = 40114
$ws.Range("L64").Value = 40114
$ws.Range("N64").Value = -40610

$ws.Range("H65").Value = 3919.5454
$ws.Range("J65").Value = 3889.375
$ws.Range("L65").Value = 19446.875
$ws.Range("N65").Value = -25686.875

$ws.Range("H67").Value = 40114
$ws.Range("J67").Value = 40114
$ws.Range("L67").Value = 40114
$ws.Range("N67").Value = -41830

$ws.Range("H70").Value = 33203.75
$ws.Range("J70").Value = 33203.75
$ws.Range("L70").Value = 33203.75
$ws.Range("N70").Value = -33833.75

$ws.Range("H73").Value = 33203.75
$ws.Range("J73").Value = 33203.75
$ws.Range("L73").Value = 33203.75
$ws.Range("N73").Value = -35387.75

$ws.Range("H81").Value = 25002088
$ws.Range("I81").Value = 1467
$ws.Range("J81").Value = 40002460
$ws.Range("K81").Value = 2934
$ws.Range("L81").Value = 80004920
$ws.Range("M81").Value = -1873
$ws.Range("N81").Value = -80007042

$ws.Range("H82").Value = 40301
$ws.Range("J82").Value = 40301
$ws.Range("L82").Value = 40301
$ws.Range("N82").Value = -41067

$ws.Range("H84").Value = 25002088
$ws.Range("I84").Value = 1467
$ws.Range("J84").Value = 40002460
$ws.Range("K84").Value = 14670
$ws.Range("L84").Value = 400024600
$ws.Range("M84").Value = -9366
$ws.Range("N84").Value = -400035208

$ws.Range("H85").Value = 40301
$ws.Range("J85").Value = 40301
$ws.Range("L85").Value = 40301
$ws.Range("N85").Value = -42953

$ws.Range("H123").Value = 27000
$ws.Range("J123").Value = 27000
$ws.Range("L123").Value = 27000
$ws.Range("N123").Value = -36800

$ws.Range("H132").Value = 2566.2917
$ws.Range("I132").Value = 1679.8182
$ws.Range("J132").Value = 3316.3845
$ws.Range("K132").Value = 5039.4546
$ws.Range("L132").Value = 9949.1535
$ws.Range("M132").Value = -2509.4546
$ws.Range("N132").Value = -15009.1535

$ws.Range("H136").Value = 5424.871
$ws.Range("I136").Value = 6177.385
$ws.Range("J136").Value = 4881.3887
$ws.Range("K136").Value = 18532.155
$ws.Range("L136").Value = 14644.1661
$ws.Range("M136").Value = -15982.155
$ws.Range("N136").Value = -19744.1661
